$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.935.04'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.845.20'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.72%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.32'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4679'
$ws.Range('E7').Value = '  +3.55%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3665'
$ws.Range('E8').Value = '  +1.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07148'
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9270'
$ws.Range('E10').Value = '  +3.58%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.938.09'
$ws.Range('E11').Value = '  +9.13%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.58'
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07704'
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.285'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.396'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.21'
$ws.Range('E16').Value = '  +3.33%  '
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008617'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.967.78'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.43'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.028'
$ws.Range('E22').Value = '  +1.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.63'
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.938'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.47'
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('E26').Value = '  +2.38%  '
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.22'
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.887'
$ws.Range('E29').Value = '  +0.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08859'
$ws.Range('E30').Value = '  +1.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.220'
$ws.Range('E31').Value = '  +2.93%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.176'
$ws.Range('E32').Value = '  +5.87%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7476'
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.801'
$ws.Range('E34').Value = '  +2.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.466'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('E36').Value = '  +1.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01940'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.970'
$ws.Range('E38').Value = '  +2.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05190'
$ws.Range('E39').Value = '  +1.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5210'
$ws.Range('E40').Value = '  +2.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.913'
$ws.Range('E41').Value = '  +2.56%  '
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.135'
$ws.Range('E43').Value = '  +0.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.49'
$ws.Range('E44').Value = '  +4.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4696'
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.008'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.59'
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.608'
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '65.03'
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06031'
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8888'
$ws.Range('E51').Value = '  +5.12%  '
